# Update "想去人数" (wanted-to-go count) figures in column F across the four
# sheets of the workbook, reflecting refreshed numbers from the data source
# (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7668
$ws.Range("F3").Value = 3611
$ws.Range("F5").Value = 3918
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 104
$ws.Range("F10").Value = 196
$ws.Range("F13").Value = 177
$ws.Range("F15").Value = 14
$ws.Range("F17").Value = 366
$ws.Range("F18").Value = 4298
$ws.Range("F19").Value = 4298
$ws.Range("F21").Value = 430
$ws.Range("F23").Value = 551
$ws.Range("F24").Value = 2668
$ws.Range("F27").Value = 3139
$ws.Range("F28").Value = 2436
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = 103
$ws.Range("F33").Value = 133
$ws.Range("F34").Value = 143
$ws.Range("F35").Value = 49
$ws.Range("F36").Value = 40
$ws.Range("F37").Value = 115
$ws.Range("F38").Value = 4587
$ws.Range("F39").Value = 565
$ws.Range("F40").Value = 345
$ws.Range("F41").Value = 64
$ws.Range("F43").Value = 890
$ws.Range("F44").Value = 290
$ws.Range("F46").Value = 1730
$ws.Range("F47").Value = 272
$ws.Range("F48").Value = 44
$ws.Range("F49").Value = 631
$ws.Range("F50").Value = 749

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 451
$ws.Range("F11").Value = 49
$ws.Range("F17").Value = 110
$ws.Range("F24").Value = 658

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 175
$ws.Range("F3").Value = 10

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 175
$ws.Range("F4").Value = 7668
$ws.Range("F5").Value = 3611
$ws.Range("F6").Value = 3918
$ws.Range("F7").Value = 104
$ws.Range("F11").Value = 196
$ws.Range("F14").Value = 177
$ws.Range("F15").Value = 14
$ws.Range("F17").Value = 366
$ws.Range("F18").Value = 4298
$ws.Range("F19").Value = 4298
$ws.Range("F21").Value = 49
$ws.Range("F23").Value = 430
$ws.Range("F24").Value = 551
$ws.Range("F25").Value = 2668
$ws.Range("F28").Value = 2436
$ws.Range("F31").Value = 133
$ws.Range("F32").Value = 143
$ws.Range("F33").Value = 49
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 115
$ws.Range("F36").Value = 110
$ws.Range("F38").Value = 4587
$ws.Range("F40").Value = 565
$ws.Range("F42").Value = 64
$ws.Range("F44").Value = 890
$ws.Range("F45").Value = 290
$ws.Range("F46").Value = 1730
$ws.Range("F47").Value = 272
$ws.Range("F48").Value = 44
$ws.Range("F49").Value = 631
$ws.Range("F50").Value = 749
